$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("H1").Value = 78.194545454545462
$ws.Range("H2").Value = 69.400000000000006
$ws.Range("H3").Value = 55.017272727272719
